$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Workbook-level view tweak: narrower window width.
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Width = 500

# ---------------------------------------------------------------------------
# 2) "Service Contacts" sheet: move delivery_organisation_path (old column R)
#    to be right before practitioner_key (old column D), shifting the
#    columns in between one place to the right.
# ---------------------------------------------------------------------------
$ws = $wb.Sheets.Item("Service Contacts")

for ($r = 1; $r -le 3; $r++) {
    # Save the value that currently lives in column R (18) - this is
    # delivery_organisation_path - before it gets overwritten by the shift
    # below.
    $movedCell = $ws.Cells.Item($r, 18)
    $movedValue = $movedCell.Value2

    # Shift columns D..Q (4..17) right into E..R (5..18), working from the
    # right so we never clobber a value before it has been read. Every
    # shifted destination is reset to the default "Normal" style first so no
    # stale formatting survives the move; it is re-applied below wherever
    # the moved-from cell actually had it.
    for ($c = 18; $c -ge 5; $c--) {
        $src = $ws.Cells.Item($r, $c - 1)
        $dst = $ws.Cells.Item($r, $c)
        $dst.Style = "Normal"
        $dst.Value2 = $src.Value2
    }

    # Drop the saved delivery_organisation_path value into the now-empty
    # column D.
    $dstD = $ws.Cells.Item($r, 4)
    $dstD.Style = "Normal"
    $dstD.Value2 = $movedValue
}

# Row 3's practitioner_key cell (old R3/PHN999:NFP01) carried an explicit
# black-font style; it now lives in D3.
$ws.Cells.Item(3, 4).Font.Color = 0

# Column widths: the old width="19" formatting (columns P:R) now applies to
# the new column D and to columns Q:R.
$ws.Columns.Item(4).ColumnWidth = 18.14
$ws.Columns.Item(17).ColumnWidth = 18.14
$ws.Columns.Item(18).ColumnWidth = 18.14

# View: selection on D1 (whole column), no frozen/scrolled top-left cell.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D1:D1048576").Select()

# ---------------------------------------------------------------------------
# 3) View-only selection tweaks on the other sheets (no data changes).
# ---------------------------------------------------------------------------
$orgs = $wb.Sheets.Item("Organisations")
$orgs.Range("H1:J3").Select()

$k10 = $wb.Sheets.Item("K10+")
$k10.Range("F1:F5").Select()

$k5 = $wb.Sheets.Item("K5")
$k5.Range("F1:F5").Select()

$sdq = $wb.Sheets.Item("SDQ")
$sdq.Application.ActiveWindow.ScrollColumn = 1
$sdq.Range("F4").Select()

$ws.Select()
